$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (Changed) date column C was stamped with a new date,
# shifting every row's value from 45180 to 45181 (one day later) for
# rows 2 through 27.
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value = 45181
    }
}
